# Applies the checklist update for Test 14 (countup_pause_countup_pclk2.v)
# and Test 15 (countdw_pause_countdw_pclk2.v): corrects the swapped
# pass/faulty wording in the last sentence of each Description/Requirement,
# marks both rows as "Pass" (matching the styling already used by other
# Pass rows), and records the completion date "Oct 11th" in Start/End Date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24: Test 14 (countup_pause_countup_pclk2.v) ---
$e24 = @'
-Enter a random number less than 255.
-The timer is enabled and counts up with an internal clock of pclk × 2, starting from the number read above.
-It runs for a while, then pauses for a duration of COUNT_PAUSE time units. 
-After the pause, it checks if an overflow occurred. If an overflow is detected, it displays a "faulty" message; otherwise, it displays a "pass" (normal operation).
-Keep the operation condition unchanged for the timer to count up and set the EN bit to 1. The timer continues counting up and, upon reaching 255, transfes to 0.
-Check if the overflow status is triggered. If so, display a "faulty" message (normal operation); otherwise, display a "pass" message.
'@
$ws.Range("E24").Value = $e24

$f24 = @'
-Write a random value to the TDR (address = 0).
-Load the value from TDR into the TCNT register.
-Set the conditions for operation, including disabling the LOAD bit, setting the count-up bit, configuring the internal clock to be equivalent to pclk × 2, and finally enabling the EN bit to put the timer into operation.
-The timer runs and counts up for a while, then disable the EN bit to stop the timer during the PAUSE time units.
-After the pause, check the overflow status and display the appropriate message: indicate a fault if overflow is triggered, or normal operation if overflow has not yet been triggered.
-Enable the timer (set the EN bit to 1). The timer continues counting up from where it stopped. Once the timer exceeds 255, check if the overflow status is triggered. If so, display a "faulty" message (normal operation); otherwise, display a "pass" message.
'@
$ws.Range("F24").Value = $f24

$g24 = $ws.Range("G24")
$g24.Value = "Pass"
$g24.Font.Color = $ws.Range("G9").Font.Color

$ws.Range("H24").Value = "Oct 11th"
$ws.Range("I24").Value = "Oct 11th"

# --- Row 25: Test 15 (countdw_pause_countdw_pclk2.v) ---
$e25 = @'
Enter a random number less than 255.
-The timer is enabled and counts down with an internal clock of pclk × 2, starting from the number read above.
-It runs for a while, then pauses for a duration of COUNT_PAUSE time units.
- After the pause, it checks if an underflow occurred. If an underflow is detected, it displays a "faulty" message; otherwise, it displays a "pass" (normal operation).
- Keep the operation condition for the timer as before: counting-down, pclk2 and set the EN bit to 1. 
-The timer continues counting down and, upon reaching below 0, converts to 255.
-Check if the underflow status is triggered. If so, display a "faulty" message (normal operation); otherwise, display a "pass" message.
'@
$ws.Range("E25").Value = $e25

$f25 = @'
-Write a random value to the TDR (address = 0).
-Load the value from TDR into the TCNT register.
-Set the conditions for operation, including disabling the LOAD bit, setting the count-down bit, configuring the internal clock to be equivalent to pclk × 2, and finally enabling the EN bit to put the timer into operation.
-The timer runs and counts down for a while, then disable the EN bit to stop the timer during the PAUSE time units.
-After the pause, check the underflow status and display the appropriate message: indicate a fault if underflow is triggered, or normal operation if underflow has not yet been triggered.
-Enable the timer (set the EN bit to 1). The timer continues counting down from where it stopped. Once the timer goes below 0, transfers to 255, check if the underflow status is triggered. If so, display a "faulty" message (normal operation); otherwise, display a "pass" message.
'@
$ws.Range("F25").Value = $f25

$g25 = $ws.Range("G25")
$g25.Value = "Pass"
$g25.Font.Color = $ws.Range("G9").Font.Color

$ws.Range("H25").Value = "Oct 11th"
$ws.Range("I25").Value = "Oct 11th"

# --- Update the saved view state: scroll/selection now centred on row 25 ---
$ws.Range("G25:I25").Select()
$excel.ActiveWindow.Zoom = 100
